# Insert a new record row at row 303, shifting the existing rows
# 303:348 down to 304:349 (Excel "insert, shift cells down" behaviour).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(303).Insert()

# Populate the newly inserted row with the new data point.
$ws.Range("A303").Value = 1
$ws.Range("B303").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C303").Value = "Arica y Parinacota"
$ws.Range("D303").Value = 44951
$ws.Range("E303").Value = 15
$ws.Range("F303").Value = "Fruta"
$ws.Range("G303").Value = 100108
$ws.Range("H303").Value = "Tropicales y subtropicales"
$ws.Range("I303").Value = 100108006
$ws.Range("J303").Value = "Plátano"
$ws.Range("K303").Value = "Barraganete"
$ws.Range("L303").Value = "Primera"
$ws.Range("M303").Value = 300
$ws.Range("N303").Value = 24000
$ws.Range("O303").Value = 25000
$ws.Range("P303").Value = 24667
$ws.Range("Q303").Value = "$/caja 20 kilos"
$ws.Range("R303").Value = "Ecuador"
$ws.Range("S303").Value = 1233
$ws.Range("T303").Value = 20
